# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Wed Jul 12 04:39:29 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is @(CellRef, NewValue). Price-column values that look like plain
# numbers are prefixed with a leading apostrophe so Excel stores them as text
# (matching the workbook convention of keeping Price/Volume as literal strings)
# instead of silently parsing them into numeric cells.
$updates = @(
  @('D2', '30.610.25'),
  @('D3', '1.881.07'),
  @('E3', '  +0.09%  '),
  @('D4', '''0.9997'),
  @('E4', '  -0.04%  '),
  @('D5', '''249.64'),
  @('E5', '  +1.20%  '),
  @('D6', '''0.9999'),
  @('E6', '  -0.03%  '),
  @('D7', '''0.4753'),
  @('E7', '  -0.12%  '),
  @('D8', '''0.2936'),
  @('E8', '  +1.15%  '),
  @('E9', '  +0.04%  '),
  @('E10', '  +0.27%  '),
  @('D11', '''0.07751'),
  @('E11', '  +0.11%  '),
  @('D12', '''96.77'),
  @('E12', '  -0.09%  '),
  @('D13', '''0.7382'),
  @('E13', '  -0.19%  '),
  @('D14', '1.880.06'),
  @('E14', '  +0.04%  '),
  @('D15', '''5.232'),
  @('E15', '  +2.02%  '),
  @('D16', '''274.42'),
  @('E16', '  +0.57%  '),
  @('D17', '30.702.00'),
  @('E17', '  +0.75%  '),
  @('E18', '  -3.43%  '),
  @('B19', 'Dai'),
  @('C19', 'https://coinranking.com/coin/MoTuySvg7+dai-dai'),
  @('D19', '''1.0000'),
  @('E19', '  -0.02%  '),
  @('B20', 'ShibaInu'),
  @('C20', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'),
  @('D20', '''0.000007525'),
  @('E20', '  -0.81%  '),
  @('D21', '2.127.48'),
  @('E21', '  +0.02%  '),
  @('D22', '''5.333'),
  @('E22', '  +1.98%  '),
  @('E23', '  -0.01%  '),
  @('D24', '''6.232'),
  @('E24', '  +0.96%  '),
  @('D25', '''9.221'),
  @('E25', '  -0.91%  '),
  @('D26', '''163.68'),
  @('E27', '  -0.32%  '),
  @('D28', '''1.905'),
  @('E28', '  -2.00%  '),
  @('E29', '  -2.05%  '),
  @('D30', '''0.09701'),
  @('E30', '  -2.78%  '),
  @('D31', '''1.506'),
  @('E31', '  -0.45%  '),
  @('D32', '''4.277'),
  @('E32', '  -0.97%  '),
  @('D33', '''4.148'),
  @('E33', '  +2.12%  '),
  @('D34', '''0.04862'),
  @('E34', '  +1.74%  '),
  @('D35', '''1.124'),
  @('E35', '  -0.15%  '),
  @('D36', '''0.6980'),
  @('E36', '  -0.21%  '),
  @('D37', '''2.718'),
  @('E38', '  +1.94%  '),
  @('E39', '  +1.59%  '),
  @('D40', '''6.293'),
  @('E40', '  -0.87%  '),
  @('D41', '''74.91'),
  @('E41', '  +6.89%  '),
  @('D42', '''2.025'),
  @('E42', '  +4.44%  '),
  @('D43', '''0.4243'),
  @('E43', '  +1.56%  '),
  @('D44', '''0.8406'),
  @('E44', '  +0.16%  '),
  @('D45', '''0.9996'),
  @('E45', '  -0.02%  '),
  @('D46', '''102.49'),
  @('E46', '  -0.21%  '),
  @('D47', '''9.360'),
  @('E47', '  +0.41%  '),
  @('D48', '''7.039'),
  @('E48', '  -0.70%  '),
  @('D49', '''35.59'),
  @('E49', '  +0.24%  '),
  @('D50', '''914.95'),
  @('E50', '  -0.25%  '),
  @('D51', '''0.05748'),
  @('E51', '  +2.07%  ')
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
